{"js": "// The template held an M2Doc field `{m:'doc.html'.fromHTMLURI()}` encoded as a\n// real Word field (fldChar begin/end + instrText runs) around a `_GoBack`\n// bookmark. The parser was updated (TokenIteratorFieldRewriterSplit) to work\n// off literal `{`...`}` delimited text runs instead, so we rewrite that field\n// in place into plain `<w:t>` runs carrying the same token pieces, keeping the\n// bookmark anchored between the \"doc.html\" and \"'.fromHTMLURI()\" pieces, and\n// dropping the now-unused field character machinery entirely.\n\nconst fields = context.document.body.fields;\nfields.load(\"items\");\nawait context.sync();\n\n// Locate the M2Doc field (its code is the script text between the braces we\n// are about to materialize). There is exactly one in this template, but\n// search defensively in case more fields are ever added around it.\nlet target = null;\nfor (const f of fields.items) {\n  f.load(\"code\");\n}\nawait context.sync();\nfor (const f of fields.items) {\n  if (f.code && f.code.indexOf(\"fromHTMLURI\") !== -1) {\n    target = f;\n    break;\n  }\n}\nif (!target && fields.items.length > 0) {\n  target = fields.items[0];\n}\n\n// The field's Result range sits inside the same (single) paragraph as the\n// field itself, so grabbing its paragraph gives us the exact paragraph to\n// rewrite without guessing indices.\nconst paras = target.result.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nconst fieldParagraph = paras.items[0];\n\n// Rebuild the paragraph as literal text runs: \"{\", \"m\", \":\", \"'\", \"doc.html\",\n// <bookmark untouched>, \"'.fromHTMLURI()\", \"}\". This mirrors exactly how the\n// instrText runs were split before, only now as <w:t> runs with no field\n// delimiters, and the bookmark keeps its original position in the run\n// sequence.\nconst replacementOoxml =\n  '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p w:rsidR=\"00C52979\" w:rsidRDefault=\"00C52979\" w:rsidP=\"00F5495F\">' +\n  '<w:r><w:t>{</w:t></w:r>' +\n  '<w:r><w:t>m</w:t></w:r>' +\n  '<w:r><w:t>:</w:t></w:r>' +\n  \"<w:r><w:t>'</w:t></w:r>\" +\n  '<w:r><w:t>doc.html</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  \"<w:r><w:t>'.fromHTMLURI()</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\n\nfieldParagraph.insertOoxml(replacementOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The template held an M2Doc field `{m:'doc.html'.fromHTMLURI()}` encoded as a\n# real Word field (fldChar begin/end + instrText runs) wrapping a `_GoBack`\n# bookmark. The parser was updated (TokenIteratorFieldRewriterSplit) to work\n# off literal `{`...`}` delimited text instead, so we rewrite that field in\n# place into plain text carrying the same token content, re-creating the\n# bookmark at the same spot in the token stream (right after \"doc.html\", the\n# same place it sat between the old instrText runs) and dropping the field\n# character machinery (fldChar begin/end + the padding instrText spaces)\n# entirely.\n\n$d = $word.ActiveDocument\n\n# Locate the M2Doc field. There is exactly one in this template; search\n# defensively by its code content in case more fields are ever added.\n$targetField = $null\nfor ($i = 1; $i -le $d.Fields.Count; $i++) {\n    $candidate = $d.Fields.Item($i)\n    if ($candidate.Code.Text -like \"*fromHTMLURI*\") {\n        $targetField = $candidate\n        break\n    }\n}\nif ($targetField -eq $null -and $d.Fields.Count -gt 0) {\n    $targetField = $d.Fields.Item(1)\n}\n\n# Pull the literal script text out of the field code, dropping the single\n# leading/trailing padding space that used to live in their own instrText\n# runs (those runs are simply gone now, not replaced by anything).\n$code = $targetField.Code.Text.Trim()\n\n# The bookmark used to sit right after the \"doc.html\" instrText run and\n# before the \"'.fromHTMLURI()\" instrText run; keep it in that same spot in\n# the token stream.\n$marker = \"doc.html\"\n$splitAt = $code.IndexOf($marker) + $marker.Length\n$beforeBookmark = $code.Substring(0, $splitAt)\n$afterBookmark = $code.Substring($splitAt)\n\n# Remember where the field's paragraph begins so we can reinsert plain text\n# at the exact same spot once the field (and its bookmark) are gone. (Code /\n# Result ranges don't reliably scope their own .Paragraphs collection here,\n# so find the owning paragraph by checking which one reports the field.)\n$insertStart = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $candidateParagraph = $d.Paragraphs.Item($i)\n    if ($candidateParagraph.Range.Fields.Count -gt 0) {\n        $insertStart = $candidateParagraph.Range.Start\n        break\n    }\n}\n\n$targetField.Delete()\n\n# Insert the whole literal replacement as a single, stable block of text\n# first. Only once that text exists do we splice the bookmark back in at its\n# spot inside it - adding a bookmark into the middle of settled text reliably\n# splits the run in two around it; doing the insertions in the other order\n# (text-then-more-text around an already-placed bookmark) does not rebase the\n# bookmark the way real Word would, so the ordering here matters.\n$literalText = \"{\" + $beforeBookmark + $afterBookmark + \"}\"\n$d.Range($insertStart, $insertStart).Text = $literalText\n\n$bookmarkPos = $insertStart + 1 + $beforeBookmark.Length\n$d.Bookmarks.Add(\"_GoBack\", $d.Range($bookmarkPos, $bookmarkPos))\n"}
